$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 24.06000000000032
$ws.Range("H2").Value = 0.0000000000000001801578944625001
$ws.Range("K2").Value = 45.34990817767195
$ws.Range("L2").Value = "[35.813706615624525, 54.88610973971937]"
$ws.Range("O2").Value = 1.817658211986888
$ws.Range("P2").Value = "[1.591237119836273, 2.0440793041375036]"
$ws.Range("S2").Value = 64.00163700583423
$ws.Range("T2").Value = "[58.577744732761346, 69.4255292789071]"
$ws.Range("W2").Value = 17.09969969969993
$ws.Range("X2").Value = 16.23267267267289
$ws.Range("Y2").Value = 17.96672672672696
$ws.Range("E3").Value = 23.88000000000029
$ws.Range("G3").Value = 0.0000000000000001110223024625157
$ws.Range("H3").Value = 0.0000000000000001801578944625001
$ws.Range("K3").Value = 53.04156095954623
$ws.Range("L3").Value = "[42.90418885454193, 63.178933064550534]"
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1.528342372016656
$ws.Range("P3").Value = "[1.3145002294299628, 1.7421845146033492]"
$ws.Range("S3").Value = 62.11733638749115
$ws.Range("T3").Value = "[55.512810369681944, 68.72186240530036]"
$ws.Range("W3").Value = 18.07135135135157
$ws.Range("X3").Value = 17.25861861861883
$ws.Range("Y3").Value = 18.88408408408432
